$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sprint Backlog Burndown: record 0.5 remaining for "Implement storing
# location information" in Week 1 (cell D3). The SUM formula in D27 (and
# the dependent burndown chart series) picks this up automatically.
$ws.Range("D3").Value = 0.5

# Leave the sheet with D4 as the active selection, matching the saved file.
$ws.Range("D4").Select()

$wb.Save()
